# Apply edits to statement_113.xlsx per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3: account holder name / card number change
$ws.Range("C2").Value = "Hartmut"
# Leading apostrophe forces this long numeric-looking card number to stay
# text (otherwise it gets auto-coerced into a double and loses precision).
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Row 5: opening balance label/date
$ws.Range("D5").Value = "KONTOSTAND AM 28.11.2024"

# Row 6
$ws.Range("B6").Value = "29.11."
$ws.Range("C6").Value = "30.11."
$ws.Range("D6").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 7969566"
$ws.Range("E6").Value = "83,02-"

# Row 7
$ws.Range("B7").Value = "03.12."
$ws.Range("C7").Value = "04.12."
$ws.Range("D7").Value = "ZALANDO MKTPLC EU OLFPBX"
$ws.Range("E7").Value = "137,70-"

# Row 8
$ws.Range("B8").Value = "04.12."
$ws.Range("C8").Value = "05.12."
$ws.Range("D8").Value = "RECHNUNG VODAFONE GMBH 84471865"
$ws.Range("E8").Value = "38,04-"

# Row 9
$ws.Range("B9").Value = "06.12."
$ws.Range("C9").Value = "07.12."
$ws.Range("D9").Value = "KARTENZAHLUNG JET TANKSTELLE"
$ws.Range("E9").Value = "40,92-"

# Row 10
$ws.Range("B10").Value = "10.12."
$ws.Range("C10").Value = "11.12."
$ws.Range("D10").Value = "BEITRAG Allianz SE K-33842487"
$ws.Range("E10").Value = "57,37-"

# Row 11: transaction removed, cells cleared (E11 reverts to the plain
# right-aligned+centered+wrap style instead of the simple right-aligned one)
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("E11").VerticalAlignment = -4108
$ws.Range("E11").WrapText = $true

# Row 12: closing balance label/date + amount
$ws.Range("D12").Value = "KONTOSTAND AM 12.12.2024"
$ws.Range("E12").Value = "357,05-"

# Row 13: next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 18.12.2024"
